# ---------------------------------------------------------------------------
# additionIn50.docx  (0.3.1 + README changed)
#
# The document starts out as a single empty paragraph. This turns it into a
# "64 addition facts" timed drill sheet:
#   1. A "date / time / score" header line.
#   2. A 4-column table: a blank spacer row, then 16 rows of two-digit
#      addition problems (64 problems total).
#   3. The Normal style picks up exact 29pt (580 twips) line spacing so the
#      printed rows line up with handwriting space.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# The 16 rows x 4 columns of addition facts that make up the drill table.
$facts = @(
    @("16+15=", "13+12=", "35+16=", "22+19="),
    @("38+18=", "18+12=", "39+37=", "42+17="),
    @("33+16=", "32+16=", "45+18=", "36+24="),
    @("43+36=", "22+14=", "48+33=", "33+17="),
    @("22+13=", "39+14=", "34+14=", "45+27="),
    @("26+17=", "32+19=", "18+16=", "14+13="),
    @("35+23=", "47+12=", "25+14=", "36+17="),
    @("38+32=", "27+14=", "33+24=", "42+39="),
    @("28+16=", "26+13=", "23+12=", "46+38="),
    @("14+12=", "27+19=", "36+23=", "28+14="),
    @("32+12=", "48+32=", "18+13=", "49+48="),
    @("18+14=", "39+26=", "23+13=", "43+42="),
    @("37+15=", "18+15=", "47+45=", "34+29="),
    @("18+17=", "48+18=", "49+46=", "33+28="),
    @("29+16=", "38+24=", "43+25=", "27+22="),
    @("32+14=", "25+23=", "43+32=", "44+26=")
)

$wNs = "xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`""

# Header paragraph: "date: ___ time: ___ ___/64"
$headerText = "date:                               time:                               ___/64"
$bodyXml = "<w:p $wNs><w:r><w:t>$headerText</w:t></w:r></w:p>"

# Table properties: autofit layout, tblLook banding flags -> val 04A0.
$bodyXml += "<w:tbl $wNs>"
$bodyXml += "<w:tblPr><w:tblW w:type=`"auto`" w:w=`"0`"/><w:tblLayout w:type=`"autofit`"/>" +
            "<w:tblLook w:firstColumn=`"1`" w:firstRow=`"1`" w:lastColumn=`"0`" w:lastRow=`"0`" " +
            "w:noHBand=`"0`" w:noVBand=`"1`" w:val=`"04A0`"/></w:tblPr>"
$bodyXml += "<w:tblGrid><w:gridCol w:w=`"2256`"/><w:gridCol w:w=`"2256`"/><w:gridCol w:w=`"2256`"/><w:gridCol w:w=`"2256`"/></w:tblGrid>"

# Leading blank spacer row, columns sized 2256 each.
$bodyXml += "<w:tr>"
for ($i = 0; $i -lt 4; $i++) {
    $bodyXml += "<w:tc><w:tcPr><w:tcW w:type=`"dxa`" w:w=`"2256`"/></w:tcPr><w:p/></w:tc>"
}
$bodyXml += "</w:tr>"

# The 16 rows of addition facts, columns sized 4209 each.
foreach ($row in $facts) {
    $bodyXml += "<w:tr>"
    foreach ($cell in $row) {
        $bodyXml += "<w:tc><w:tcPr><w:tcW w:type=`"dxa`" w:w=`"4209`"/></w:tcPr>" +
                    "<w:p><w:r><w:t>$cell</w:t></w:r></w:p></w:tc>"
    }
    $bodyXml += "</w:tr>"
}
$bodyXml += "</w:tbl>"

# Replace the document's entire current content (the one empty paragraph,
# including its trailing paragraph mark) with the header + table so nothing
# empty is left behind before the final sectPr.
$fullRange = $d.Range(0, $d.Content.End)
$fullRange.InsertXML($bodyXml)

# Normal style: exact 29pt (580 twips) line spacing.
$normalStyle = $d.Styles("Normal")
$normalStyle.ParagraphFormat.LineSpacingRule = 4   # wdLineSpaceExactly
$normalStyle.ParagraphFormat.LineSpacing = 29       # points (29pt = 580 twips)
